$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new Price text, new Volume(1h) percent text)
$updates = @{
  2 = @("22.447.47", "-0.04%")
  3 = @("1.567.79", "-0.31%")
  4 = @("1.003", "+0.11%")
  5 = @("1.002", "+0.10%")
  6 = @("288.45", "-0.60%")
  7 = @("0.3729", "+0.83%")
  8 = @("48.24", "-3.28%")
  9 = @("0.3321", "-1.97%")
  10 = @("0.07477", "-1.20%")
  11 = @("1.131", "-1.72%")
  12 = @("1.003", "+0.09%")
  13 = @("20.77", "-2.30%")
  14 = @("5.962", "-1.23%")
  15 = @("6.906", "-1.27%")
  16 = @("1.563.37", "-1.42%")
  17 = @("0.00001114", "-0.89%")
  18 = @("0.06769", "-0.34%")
  19 = @("87.94", "-2.94%")
  20 = @("1.003", "+0.01%")
  21 = @("6.355", "-0.15%")
  22 = @("16.41", "-0.15%")
  23 = @("12.08", "-1.02%")
  24 = @("22.455.37", "-0.07%")
  25 = @("2.391", "+0.88%")
  26 = @("2.564", "-4.14%")
  27 = @("154.11", "+3.07%")
  28 = @("19.64", "-2.01%")
  29 = @("5.014", "-0.81%")
  30 = @("124.04", "-0.88%")
  31 = @("1.744.86", "-0.22%")
  32 = @("1.054", "-1.17%")
  33 = @("2.014", "-0.19%")
  34 = @("6.125", "-1.55%")
  35 = @("9.626", "-2.17%")
  36 = @("0.08284", "-1.30%")
  37 = @("0.02457", "-0.96%")
  38 = @("0.2275", "-1.30%")
  39 = @("0.06384", "-2.81%")
  40 = @("5.359", "-1.50%")
  41 = @("1.287", "-4.59%")
  42 = @("0.6282", "+0.40%")
  43 = @("11.24", "-0.83%")
  44 = @("1.003", "+0.13%")
  45 = @("13.88", "-1.34%")
  46 = @("0.6130", "+4.25%")
  47 = @("3.777", "-0.67%")
  48 = @("2.045", "-1.48%")
  49 = @("125.57", "-1.91%")
  50 = @("1.212", "-2.44%")
  51 = @("0.07241", "-0.88%")
}

foreach ($row in $updates.Keys) {
  $priceText = $updates[$row][0]
  $volText   = $updates[$row][1]

  $dCell = $ws.Range("D$row")
  # Values such as '1.003' parse as numbers in Excel; the source data stores
  # them as plain text, so force text formatting whenever the string has at
  # most one '.' (i.e. could otherwise be auto-converted to a number).
  $dotCount = ($priceText.ToCharArray() | Where-Object { $_ -eq '.' }).Count
  if ($dotCount -le 1) {
    $dCell.NumberFormat = "@"
    $dCell.Value = $priceText
    $dCell.Style = "Normal"
  } else {
    $dCell.Value = $priceText
  }

  $ws.Range("E$row").Value = "  $volText  "
}
